$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells are treated as text so numeric-looking strings
# (e.g. "0.06630", "73.00", "1.009") keep their exact literal formatting
# instead of being coerced into Excel numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "21.375.32"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.28%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.548.60"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +4.94%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.9680"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.48%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "282.30"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3639"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3204"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +5.26%  "
$ws.Range("B9").NumberFormat = "@"
$ws.Range("B9").Value = "Polygon"
$ws.Range("C9").NumberFormat = "@"
$ws.Range("C9").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "1.117"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +6.53%  "
$ws.Range("B10").NumberFormat = "@"
$ws.Range("B10").Value = "OKB"
$ws.Range("C10").NumberFormat = "@"
$ws.Range("C10").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "40.95"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +3.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.06951"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +5.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.07%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "19.07"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +5.17%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.728"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +5.05%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.429"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +4.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001058"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +2.88%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9692"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.03%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.547.21"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +4.64%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06144"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +4.23%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "73.00"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +5.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.756"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +5.41%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +6.49%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "11.48"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +4.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.318"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +3.13%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "21.364.60"
$ws.Range("B26").NumberFormat = "@"
$ws.Range("B26").Value = "LidoDAOToken"
$ws.Range("C26").NumberFormat = "@"
$ws.Range("C26").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.313"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +8.75%  "
$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = "Monero"
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "147.69"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +4.40%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "17.87"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  +4.23%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.718.74"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +5.27%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "119.14"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  +5.14%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.037"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +3.46%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.8943"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  +10.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.265"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +6.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08082"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +2.75%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.559"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.228"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.88%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.017"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +6.18%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05920"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.2009"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +7.24%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.02141"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +5.07%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "10.89"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  +4.84%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "7.952"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.73%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.9689"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +0.62%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.5541"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +5.25%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "12.57"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +4.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.589"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.42%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5522"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +7.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "122.54"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.05%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.889"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +6.67%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.06630"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +2.84%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "70.34"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +4.95%  "
